$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 762 (shifts existing rows 762-790 down to 766-794)
$ws.Rows("762:765").Insert()

# Shared/constant column values for this product block
$mercadoId  = 9
$mercado    = 'Vega Central Mapocho de Santiago'
$region     = 'Metropolitana'
$codreg     = 13
$tipo       = 'Fruta'
$productoId = 100108
$producto   = 'Tropicales y subtropicales'
$categoriaId= 100108005
$categoria  = 'Piña'
$variedad   = 'Caramelo'
$origen     = 'Ecuador'
$fecha      = 44509

# Data specific to each new quality row (Especial / Primera / Segunda / Tercera)
$rows = @(
    @{ Row=762; Calidad='Especial'; Volumen=35; PMin=17000; PMax=18000; PProm=17429; Unidad='$/caja 10 unidades'; PrecioKg=1743; KgUnidad=10 },
    @{ Row=763; Calidad='Primera';  Volumen=30; PMin=17000; PMax=18000; PProm=17333; Unidad='$/caja 12 unidades'; PrecioKg=1444; KgUnidad=12 },
    @{ Row=764; Calidad='Segunda';  Volumen=45; PMin=17000; PMax=18000; PProm=17556; Unidad='$/caja 14 unidades'; PrecioKg=1254; KgUnidad=14 },
    @{ Row=765; Calidad='Tercera';  Volumen=25; PMin=17000; PMax=18000; PProm=17600; Unidad='$/caja 16 unidades'; PrecioKg=1100; KgUnidad=16 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $mercadoId
    $ws.Cells.Item($rowNum, 2).Value = $mercado
    $ws.Cells.Item($rowNum, 3).Value = $region
    $ws.Cells.Item($rowNum, 4).Value = $fecha
    $ws.Cells.Item($rowNum, 5).Value = $codreg
    $ws.Cells.Item($rowNum, 6).Value = $tipo
    $ws.Cells.Item($rowNum, 7).Value = $productoId
    $ws.Cells.Item($rowNum, 8).Value = $producto
    $ws.Cells.Item($rowNum, 9).Value = $categoriaId
    $ws.Cells.Item($rowNum, 10).Value = $categoria
    $ws.Cells.Item($rowNum, 11).Value = $variedad
    $ws.Cells.Item($rowNum, 12).Value = $r.Calidad
    $ws.Cells.Item($rowNum, 13).Value = $r.Volumen
    $ws.Cells.Item($rowNum, 14).Value = $r.PMin
    $ws.Cells.Item($rowNum, 15).Value = $r.PMax
    $ws.Cells.Item($rowNum, 16).Value = $r.PProm
    $ws.Cells.Item($rowNum, 17).Value = $r.Unidad
    $ws.Cells.Item($rowNum, 18).Value = $origen
    $ws.Cells.Item($rowNum, 19).Value = $r.PrecioKg
    $ws.Cells.Item($rowNum, 20).Value = $r.KgUnidad
}
